$wb = $excel.ActiveWorkbook

# --- Rename header cells on the two existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "PO Forecast"

# --- Header row ---
$newSheet.Cells.Item(1,1).Value = "ds"
$newSheet.Cells.Item(1,2).Value = "PO_Forecast"
$newSheet.Cells.Item(1,3).Value = "yhat_lower"
$newSheet.Cells.Item(1,4).Value = "yhat_upper"

# --- Data rows ---
$newSheet.Cells.Item(2,1).Value = 45067.99999999999
$newSheet.Cells.Item(2,2).Value = 34
$newSheet.Cells.Item(2,3).Value = -36.65652001541812
$newSheet.Cells.Item(2,4).Value = 107.9051231298307
$newSheet.Cells.Item(3,1).Value = 45074.99999999999
$newSheet.Cells.Item(3,2).Value = 35
$newSheet.Cells.Item(3,3).Value = -39.78278642161506
$newSheet.Cells.Item(3,4).Value = 101.9740490003794
$newSheet.Cells.Item(4,1).Value = 45102.99999999999
$newSheet.Cells.Item(4,2).Value = 36
$newSheet.Cells.Item(4,3).Value = -35.19882674923815
$newSheet.Cells.Item(4,4).Value = 105.5447469277526
$newSheet.Cells.Item(5,1).Value = 45123.99999999999
$newSheet.Cells.Item(5,2).Value = 37
$newSheet.Cells.Item(5,3).Value = -33.88185603926529
$newSheet.Cells.Item(5,4).Value = 107.5762543837228
$newSheet.Cells.Item(6,1).Value = 45130.99999999999
$newSheet.Cells.Item(6,2).Value = 38
$newSheet.Cells.Item(6,3).Value = -36.98914432990795
$newSheet.Cells.Item(6,4).Value = 103.2012740649097
$newSheet.Cells.Item(7,1).Value = 45137.99999999999
$newSheet.Cells.Item(7,2).Value = 38
$newSheet.Cells.Item(7,3).Value = -29.87817790300408
$newSheet.Cells.Item(7,4).Value = 112.9567496518735
$newSheet.Cells.Item(8,1).Value = 45158.99999999999
$newSheet.Cells.Item(8,2).Value = 39
$newSheet.Cells.Item(8,3).Value = -36.88378931795645
$newSheet.Cells.Item(8,4).Value = 112.4550521407918
$newSheet.Cells.Item(9,1).Value = 45179.99999999999
$newSheet.Cells.Item(9,2).Value = 40
$newSheet.Cells.Item(9,3).Value = -30.3254387600033
$newSheet.Cells.Item(9,4).Value = 119.40353082914
$newSheet.Cells.Item(10,1).Value = 45249.99999999999
$newSheet.Cells.Item(10,2).Value = 43
$newSheet.Cells.Item(10,3).Value = -29.82656831774446
$newSheet.Cells.Item(10,4).Value = 114.9412990370576
$newSheet.Cells.Item(11,1).Value = 45256.99999999999
$newSheet.Cells.Item(11,2).Value = 44
$newSheet.Cells.Item(11,3).Value = -26.97244124218201
$newSheet.Cells.Item(11,4).Value = 113.9993175591744
$newSheet.Cells.Item(12,1).Value = 45263.99999999999
$newSheet.Cells.Item(12,2).Value = 44
$newSheet.Cells.Item(12,3).Value = -29.44133681234654
$newSheet.Cells.Item(12,4).Value = 112.3161423256091
$newSheet.Cells.Item(13,1).Value = 45445.99999999999
$newSheet.Cells.Item(13,2).Value = 53
$newSheet.Cells.Item(13,3).Value = -16.36797189573728
$newSheet.Cells.Item(13,4).Value = 122.3490398383033
$newSheet.Cells.Item(14,1).Value = 45459.99999999999
$newSheet.Cells.Item(14,2).Value = 54
$newSheet.Cells.Item(14,3).Value = -16.05775541584238
$newSheet.Cells.Item(14,4).Value = 128.0221005051965
$newSheet.Cells.Item(15,1).Value = 45466.99999999999
$newSheet.Cells.Item(15,2).Value = 54
$newSheet.Cells.Item(15,3).Value = -16.71903872662429
$newSheet.Cells.Item(15,4).Value = 124.5531590943876
$newSheet.Cells.Item(16,1).Value = 45473.99999999999
$newSheet.Cells.Item(16,2).Value = 55
$newSheet.Cells.Item(16,3).Value = -18.39272256735205
$newSheet.Cells.Item(16,4).Value = 128.8883390622888
$newSheet.Cells.Item(17,1).Value = 45480.99999999999
$newSheet.Cells.Item(17,2).Value = 55
$newSheet.Cells.Item(17,3).Value = -19.07321450645999
$newSheet.Cells.Item(17,4).Value = 123.0166126726024
$newSheet.Cells.Item(18,1).Value = 45487.99999999999
$newSheet.Cells.Item(18,2).Value = 55
$newSheet.Cells.Item(18,3).Value = -12.00755204897514
$newSheet.Cells.Item(18,4).Value = 128.3342788198451
$newSheet.Cells.Item(19,1).Value = 45494.99999999999
$newSheet.Cells.Item(19,2).Value = 56
$newSheet.Cells.Item(19,3).Value = -19.16435501568797
$newSheet.Cells.Item(19,4).Value = 129.7366217010809
$newSheet.Cells.Item(20,1).Value = 45501.99999999999
$newSheet.Cells.Item(20,2).Value = 56
$newSheet.Cells.Item(20,3).Value = -17.19030820764165
$newSheet.Cells.Item(20,4).Value = 124.7263489072723
$newSheet.Cells.Item(21,1).Value = 45529.99999999999
$newSheet.Cells.Item(21,2).Value = 57
$newSheet.Cells.Item(21,3).Value = -12.77989327569266
$newSheet.Cells.Item(21,4).Value = 123.585379498975
$newSheet.Cells.Item(22,1).Value = 45536.99999999999
$newSheet.Cells.Item(22,2).Value = 58
$newSheet.Cells.Item(22,3).Value = -15.45431827115714
$newSheet.Cells.Item(22,4).Value = 127.8703565262863
$newSheet.Cells.Item(23,1).Value = 45550.99999999999
$newSheet.Cells.Item(23,2).Value = 58
$newSheet.Cells.Item(23,3).Value = -11.45212969559344
$newSheet.Cells.Item(23,4).Value = 132.2744813264473
$newSheet.Cells.Item(24,1).Value = 45557.99999999999
$newSheet.Cells.Item(24,2).Value = 59
$newSheet.Cells.Item(24,3).Value = -12.36536195045419
$newSheet.Cells.Item(24,4).Value = 129.6877534669947
$newSheet.Cells.Item(25,1).Value = 45564.99999999999
$newSheet.Cells.Item(25,2).Value = 59
$newSheet.Cells.Item(25,3).Value = -9.634974408542346
$newSheet.Cells.Item(25,4).Value = 132.0409078379039
$newSheet.Cells.Item(26,1).Value = 45571.99999999999
$newSheet.Cells.Item(26,2).Value = 59
$newSheet.Cells.Item(26,3).Value = -12.12008312277443
$newSheet.Cells.Item(26,4).Value = 126.9952892625436
$newSheet.Cells.Item(27,1).Value = 45578.99999999999
$newSheet.Cells.Item(27,2).Value = 60
$newSheet.Cells.Item(27,3).Value = -12.7704194812615
$newSheet.Cells.Item(27,4).Value = 131.8141472909067
$newSheet.Cells.Item(28,1).Value = 45585.99999999999
$newSheet.Cells.Item(28,2).Value = 60
$newSheet.Cells.Item(28,3).Value = -9.636474868627559
$newSheet.Cells.Item(28,4).Value = 134.5566391521223
$newSheet.Cells.Item(29,1).Value = 45599.99999999999
$newSheet.Cells.Item(29,2).Value = 61
$newSheet.Cells.Item(29,3).Value = -15.87759402567113
$newSheet.Cells.Item(29,4).Value = 130.0643884492598
$newSheet.Cells.Item(30,1).Value = 45606.99999999999
$newSheet.Cells.Item(30,2).Value = 61
$newSheet.Cells.Item(30,3).Value = -11.90805152338274
$newSheet.Cells.Item(30,4).Value = 134.5942671983327
$newSheet.Cells.Item(31,1).Value = 45627.99999999999
$newSheet.Cells.Item(31,2).Value = 62
$newSheet.Cells.Item(31,3).Value = -6.528817946115255
$newSheet.Cells.Item(31,4).Value = 137.5962996518033
$newSheet.Cells.Item(32,1).Value = 45634.99999999999
$newSheet.Cells.Item(32,2).Value = 63
$newSheet.Cells.Item(32,3).Value = -13.13833082380431
$newSheet.Cells.Item(32,4).Value = 134.5521926824343
$newSheet.Cells.Item(33,1).Value = 45641.99999999999
$newSheet.Cells.Item(33,2).Value = 63
$newSheet.Cells.Item(33,3).Value = -6.489935489280938
$newSheet.Cells.Item(33,4).Value = 128.8857067660774
$newSheet.Cells.Item(34,1).Value = 45648.99999999999
$newSheet.Cells.Item(34,2).Value = 63
$newSheet.Cells.Item(34,3).Value = -5.422625882645584
$newSheet.Cells.Item(34,4).Value = 136.4271740437584
$newSheet.Cells.Item(35,1).Value = 45655.99999999999
$newSheet.Cells.Item(35,2).Value = 64
$newSheet.Cells.Item(35,3).Value = -9.259706797637193
$newSheet.Cells.Item(35,4).Value = 133.5428492094141
$newSheet.Cells.Item(36,1).Value = 45662.99999999999
$newSheet.Cells.Item(36,2).Value = 64
$newSheet.Cells.Item(36,3).Value = -4.210749824482684
$newSheet.Cells.Item(36,4).Value = 135.4085326726402
$newSheet.Cells.Item(37,1).Value = 45669.99999999999
$newSheet.Cells.Item(37,2).Value = 64
$newSheet.Cells.Item(37,3).Value = -13.11577305438408
$newSheet.Cells.Item(37,4).Value = 133.3896481607862
$newSheet.Cells.Item(38,1).Value = 45676.99999999999
$newSheet.Cells.Item(38,2).Value = 65
$newSheet.Cells.Item(38,3).Value = -2.853186585241473
$newSheet.Cells.Item(38,4).Value = 135.4261932577126
$newSheet.Cells.Item(39,1).Value = 45683.99999999999
$newSheet.Cells.Item(39,2).Value = 65
$newSheet.Cells.Item(39,3).Value = -8.775598499105289
$newSheet.Cells.Item(39,4).Value = 137.3129284868583

# --- Formatting: copy the header style (bold/border/centered) from an
#     existing header row, and the date style from an existing date column,
#     so the new sheet re-uses the same style indices already present in the
#     workbook rather than inventing new ones. ---
$wsWeekly.Range("A1:B1").Copy() | Out-Null
$newSheet.Range("A1:D1").PasteSpecial(-4122) | Out-Null

$wsWeekly.Range("A2").Copy() | Out-Null
$newSheet.Range("A2:A39").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# Restore the originally active sheet/selection so we don't perturb
# workbook-level view state beyond what the diff specifies.
$wsWeekly.Select() | Out-Null
$wsWeekly.Range("A1").Select() | Out-Null
